$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) labels: columns F1:M1 shift left by one, and
# a new "River Stage" label is appended at M1.
$ws.Range("F1").Value = "Well ID"
$ws.Range("G1").Value = "Pumping Rate"
$ws.Range("H1").Value = "X-Coordinates"
$ws.Range("I1").Value = "Y-Coordinates"
$ws.Range("J1").Value = "Layer ID"
$ws.Range("K1").Value = "K Value"
$ws.Range("L1").Value = "D Value"
$ws.Range("M1").Value = "River Stage"

# Remove the old sample data rows 2 and 3 entirely.
$ws.Rows("2:3").Delete()
